$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 906.68134
$ws.Range("I15").Value = 906.68134
$ws.Range("K15").Value = 2720.04402
$ws.Range("M15").Value = -2551.04402
$ws.Range("H17").Value = 661.125
$ws.Range("J17").Value = 567.50616
$ws.Range("L17").Value = 1702.51848
$ws.Range("N17").Value = -2038.51848
$ws.Range("H97").Value = 6000
$ws.Range("J97").Value = 10000
$ws.Range("L97").Value = 30000
$ws.Range("N97").Value = -30992
$ws.Range("H109").Value = 29727.273
$ws.Range("J109").Value = 29727.273
$ws.Range("L109").Value = 29727.273
$ws.Range("N109").Value = -32501.273
$ws.Range("H111").Value = 1036.3334
$ws.Range("I111").Value = 129
$ws.Range("J111").Value = 1490
$ws.Range("K111").Value = 387
$ws.Range("L111").Value = 4470
$ws.Range("M111").Value = 2680
$ws.Range("N111").Value = -10604
$ws.Range("H112").Value = 1588.1455
$ws.Range("J112").Value = 1594.4445
$ws.Range("L112").Value = 4783.333500000001
$ws.Range("N112").Value = -6999.333500000001
$ws.Range("H113").Value = 6553.2354
$ws.Range("I113").Value = 4733
$ws.Range("J113").Value = 6943.2856
$ws.Range("K113").Value = 4733
$ws.Range("L113").Value = 6943.2856
$ws.Range("M113").Value = -1479
$ws.Range("N113").Value = -13451.2856
$ws.Range("H115").Value = 1094.6666
$ws.Range("I115").Value = 1094.6666
$ws.Range("J115").Value = 0
$ws.Range("K115").Value = 3283.9998
$ws.Range("L115").Value = 0
$ws.Range("M115").Value = -1716.9998
$ws.Range("N115").ClearContents()
$ws.Range("H116").Value = 509575.56
$ws.Range("I116").Value = 2003361
$ws.Range("J116").Value = 11647.066
$ws.Range("K116").Value = 2003361
$ws.Range("L116").Value = 11647.066
$ws.Range("M116").Value = -1999919
$ws.Range("N116").Value = -18531.066
$ws.Range("H117").Value = 39480
$ws.Range("J117").Value = 39480
$ws.Range("L117").Value = 39480
$ws.Range("N117").Value = -48658
$ws.Range("H118").Value = 825.8889
$ws.Range("I118").Value = 671.7778
$ws.Range("J118").Value = 980
$ws.Range("K118").Value = 2015.3334
$ws.Range("L118").Value = 2940
$ws.Range("M118").Value = -358.3334
$ws.Range("N118").Value = -6254
$ws.Range("H129").Value = 844.3913
$ws.Range("I129").Value = 370.8
$ws.Range("J129").Value = 975.94446
$ws.Range("K129").Value = 1112.4
$ws.Range("L129").Value = 2927.83338
$ws.Range("M129").Value = 3887.6
$ws.Range("N129").Value = -12927.83338
$ws.Range("H135").Value = 303.55
$ws.Range("I135").Value = 271.82144
$ws.Range("J135").Value = 377.58334
$ws.Range("K135").Value = 2446.39296
$ws.Range("L135").Value = 3398.25006
$ws.Range("M135").Value = 88.60703999999987
$ws.Range("N135").Value = -8468.25006
$ws.Range("H137").Value = 2672.5
$ws.Range("I137").Value = 1254.174
$ws.Range("J137").Value = 4225.905
$ws.Range("K137").Value = 3762.522
$ws.Range("L137").Value = 12677.715
$ws.Range("M137").Value = -1212.522
$ws.Range("N137").Value = -17777.715
$ws.Range("H138").Value = 3398.1516
$ws.Range("I138").Value = 718.27026
$ws.Range("J138").Value = 4997.4355
$ws.Range("K138").Value = 2154.81078
$ws.Range("L138").Value = 14992.3065
$ws.Range("M138").Value = 2985.18922
$ws.Range("N138").Value = -25272.3065
$ws.Range("H141").Value = 3525.0205
$ws.Range("I141").Value = 3410.372
$ws.Range("J141").Value = 4346.6665
$ws.Range("K141").Value = 10231.116
$ws.Range("L141").Value = 13039.9995
$ws.Range("M141").Value = -5051.116
$ws.Range("N141").Value = -23399.9995

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2849.883
$ws.Range("I32").Value = 2515.093
$ws.Range("K32").Value = 2515.093
$ws.Range("M32").Value = -2228.093
$ws.Range("H61").Value = 1397
$ws.Range("I61").Value = 1008.8571
$ws.Range("J61").Value = 2173.2856
$ws.Range("K61").Value = 1008.8571
$ws.Range("L61").Value = 2173.2856
$ws.Range("M61").Value = -796.8570999999999
$ws.Range("N61").Value = -2597.2856
$ws.Range("H74").Value = 3549.3242
$ws.Range("I74").Value = 3718.5356
$ws.Range("J74").Value = 3022.889
$ws.Range("K74").Value = 3718.5356
$ws.Range("L74").Value = 3022.889
$ws.Range("M74").Value = -2844.5356
$ws.Range("N74").Value = -4770.889
$ws.Range("H76").Value = 25616
$ws.Range("J76").Value = 25616
$ws.Range("L76").Value = 25616
$ws.Range("N76").Value = -26292
$ws.Range("H77").Value = 3549.3242
$ws.Range("I77").Value = 3718.5356
$ws.Range("J77").Value = 3022.889
$ws.Range("K77").Value = 18592.678
$ws.Range("L77").Value = 15114.445
$ws.Range("M77").Value = -14224.678
$ws.Range("N77").Value = -23850.445
$ws.Range("H79").Value = 25616
$ws.Range("J79").Value = 25616
$ws.Range("L79").Value = 25616
$ws.Range("N79").Value = -27956
$ws.Range("H80").Value = 40623.5
$ws.Range("J80").Value = 40623.5
$ws.Range("L80").Value = 40623.5
$ws.Range("N80").Value = -42619.5
$ws.Range("H83").Value = 40623.5
$ws.Range("J83").Value = 40623.5
$ws.Range("L83").Value = 121870.5
$ws.Range("N83").Value = -131854.5
$ws.Range("H88").Value = 16669111
$ws.Range("I88").Value = 66666664
$ws.Range("J88").Value = 3260
$ws.Range("K88").Value = 66666664
$ws.Range("L88").Value = 3260
$ws.Range("M88").Value = -66666258
$ws.Range("N88").Value = -4072
$ws.Range("H91").Value = 16669111
$ws.Range("I91").Value = 66666664
$ws.Range("J91").Value = 3260
$ws.Range("K91").Value = 66666664
$ws.Range("L91").Value = 3260
$ws.Range("M91").Value = -66665260
$ws.Range("N91").Value = -6068
$ws.Range("H107").Value = 0
$ws.Range("J107").Value = 0
$ws.Range("L107").Value = 0
$ws.Range("N107").ClearContents()
$ws.Range("H132").Value = 1891.5588
$ws.Range("I132").Value = 1226.26
$ws.Range("J132").Value = 3739.611
$ws.Range("K132").Value = 3678.78
$ws.Range("L132").Value = 11218.833
$ws.Range("M132").Value = -1148.78
$ws.Range("N132").Value = -16278.833
$ws.Range("H136").Value = 1397
$ws.Range("I136").Value = 1008.8571
$ws.Range("J136").Value = 2173.2856
$ws.Range("K136").Value = 3026.5713
$ws.Range("L136").Value = 6519.8568
$ws.Range("M136").Value = -476.5712999999996
$ws.Range("N136").Value = -11619.8568

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 686
$ws.Range("I64").Value = 678
$ws.Range("J64").Value = 690
$ws.Range("K64").Value = 678
$ws.Range("L64").Value = 690
$ws.Range("M64").Value = -453
$ws.Range("N64").Value = -1140
$ws.Range("H67").Value = 686
$ws.Range("I67").Value = 678
$ws.Range("J67").Value = 690
$ws.Range("K67").Value = 678
$ws.Range("L67").Value = 690
$ws.Range("M67").Value = 102
$ws.Range("N67").Value = -2250
$ws.Range("H94").Value = 478.85294
$ws.Range("I94").Value = 476.14816
$ws.Range("J94").Value = 489.2857
$ws.Range("K94").Value = 476.14816
$ws.Range("L94").Value = 489.2857
$ws.Range("M94").Value = -25.14816000000002
$ws.Range("N94").Value = -1391.2857

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 503.60605
$ws.Range("I22").Value = 311.25
$ws.Range("K22").Value = 311.25
$ws.Range("M22").Value = 38.75
$ws.Range("H31").Value = 8774745
$ws.Range("I31").Value = 1563.6562
$ws.Range("J31").Value = 20004416
$ws.Range("K31").Value = 1563.6562
$ws.Range("L31").Value = 20004416
$ws.Range("M31").Value = -1268.6562
$ws.Range("N31").Value = -20005006
$ws.Range("H34").Value = 8774745
$ws.Range("I34").Value = 1563.6562
$ws.Range("J34").Value = 20004416
$ws.Range("K34").Value = 1563.6562
$ws.Range("L34").Value = 20004416
$ws.Range("M34").Value = -1361.6562
$ws.Range("N34").Value = -20004820
$ws.Range("H58").Value = 2033.1805
$ws.Range("I58").Value = 1834.2727
$ws.Range("J58").Value = 2676.7058
$ws.Range("K58").Value = 1834.2727
$ws.Range("L58").Value = 2676.7058
$ws.Range("M58").Value = -1631.2727
$ws.Range("N58").Value = -3082.7058
$ws.Range("H115").Value = 30000
$ws.Range("J115").Value = 30000
$ws.Range("L115").Value = 30000
$ws.Range("N115").Value = -32350
$ws.Range("H132").Value = 2471.2432
$ws.Range("I132").Value = 1891.5483
$ws.Range("J132").Value = 5466.3335
$ws.Range("K132").Value = 5674.644899999999
$ws.Range("L132").Value = 16399.0005
$ws.Range("M132").Value = -3144.644899999999
$ws.Range("N132").Value = -21459.0005
$ws.Range("H134").Value = 4554.3423
$ws.Range("I134").Value = 6649.222
$ws.Range("J134").Value = 2668.95
$ws.Range("K134").Value = 19947.666
$ws.Range("L134").Value = 8006.849999999999
$ws.Range("M134").Value = -17412.666
$ws.Range("N134").Value = -13076.85
$ws.Range("H136").Value = 2033.1805
$ws.Range("I136").Value = 1834.2727
$ws.Range("J136").Value = 2676.7058
$ws.Range("K136").Value = 5502.8181
$ws.Range("L136").Value = 8030.117400000001
$ws.Range("M136").Value = -2952.8181
$ws.Range("N136").Value = -13130.1174

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1435.4722
$ws.Range("I5").Value = 342.72223
$ws.Range("J5").Value = 2528.2222
$ws.Range("K5").Value = 1028.16669
$ws.Range("L5").Value = 7584.6666
$ws.Range("M5").Value = -916.16669
$ws.Range("N5").Value = -7808.6666
$ws.Range("H113").Value = 517.56525
$ws.Range("I113").Value = 478.79486
$ws.Range("J113").Value = 567.9666999999999
$ws.Range("K113").Value = 1436.38458
$ws.Range("L113").Value = 1703.9001
$ws.Range("M113").Value = 733.6154199999999
$ws.Range("N113").Value = -6043.9001
$ws.Range("H131").Value = 858.29346
$ws.Range("J131").Value = 905.03613
$ws.Range("L131").Value = 2715.10839
$ws.Range("N131").Value = -12795.10839
$ws.Range("H132").Value = 2021.9286
$ws.Range("I132").Value = 856.6
$ws.Range("J132").Value = 2669.3333
$ws.Range("K132").Value = 7709.400000000001
$ws.Range("L132").Value = 24023.9997
$ws.Range("M132").Value = -5179.400000000001
$ws.Range("N132").Value = -29083.9997
$ws.Range("H135").Value = 1435.4722
$ws.Range("I135").Value = 342.72223
$ws.Range("J135").Value = 2528.2222
$ws.Range("K135").Value = 3084.50007
$ws.Range("L135").Value = 22753.9998
$ws.Range("M135").Value = -549.5000700000001
$ws.Range("N135").Value = -27823.9998

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H31").Value = 1000
$ws.Range("I31").Value = 1000
$ws.Range("J31").Value = 0
$ws.Range("K31").Value = 1000
$ws.Range("L31").Value = 0
$ws.Range("M31").Value = -708
$ws.Range("N31").ClearContents()
$ws.Range("H37").Value = 1000
$ws.Range("I37").Value = 1000
$ws.Range("J37").Value = 0
$ws.Range("K37").Value = 1000
$ws.Range("L37").Value = 0
$ws.Range("M37").Value = -723
$ws.Range("N37").ClearContents()
$ws.Range("H122").Value = 4287.2144
$ws.Range("J122").Value = 5161.6
$ws.Range("L122").Value = 15484.8
$ws.Range("N122").Value = -20384.8
$ws.Range("H123").Value = 10861.053
$ws.Range("J123").Value = 10861.053
$ws.Range("L123").Value = 10861.053
$ws.Range("N123").Value = -15761.053
$ws.Range("H130").Value = 0
$ws.Range("J130").Value = 0
$ws.Range("L130").Value = 0
$ws.Range("N130").ClearContents()
$ws.Range("H132").Value = 2512.1628
$ws.Range("I132").Value = 1556.2222
$ws.Range("J132").Value = 4125.3125
$ws.Range("K132").Value = 4668.6666
$ws.Range("L132").Value = 12375.9375
$ws.Range("M132").Value = -2138.6666
$ws.Range("N132").Value = -17435.9375
$ws.Range("H140").Value = 38746.9
$ws.Range("J140").Value = 38746.9
$ws.Range("L140").Value = 38746.9
$ws.Range("N140").Value = -49106.9

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 33857.582
$ws.Range("I22").Value = 112275.78
$ws.Range("J22").Value = 1777.409
$ws.Range("K22").Value = 112275.78
$ws.Range("L22").Value = 1777.409
$ws.Range("M22").Value = -111980.78
$ws.Range("N22").Value = -2367.409
$ws.Range("H27").Value = 33857.582
$ws.Range("I27").Value = 112275.78
$ws.Range("J27").Value = 1777.409
$ws.Range("K27").Value = 112275.78
$ws.Range("L27").Value = 1777.409
$ws.Range("M27").Value = -112168.78
$ws.Range("N27").Value = -1991.409
$ws.Range("H69").Value = 0
$ws.Range("J69").Value = 0
$ws.Range("L69").Value = 0
$ws.Range("N69").ClearContents()
$ws.Range("H72").Value = 0
$ws.Range("J72").Value = 0
$ws.Range("L72").Value = 0
$ws.Range("N72").ClearContents()
$ws.Range("H74").Value = 39000
$ws.Range("J74").Value = 50000
$ws.Range("L74").Value = 50000
$ws.Range("N74").Value = -51996
$ws.Range("H77").Value = 39000
$ws.Range("J77").Value = 50000
$ws.Range("L77").Value = 150000
$ws.Range("N77").Value = -159984
$ws.Range("H119").Value = 35000
$ws.Range("J119").Value = 35000
$ws.Range("L119").Value = 35000
$ws.Range("N119").Value = -44676
$ws.Range("H120").Value = 0
$ws.Range("J120").Value = 0
$ws.Range("L120").Value = 0
$ws.Range("N120").ClearContents()
$ws.Range("H132").Value = 9118.66
$ws.Range("I132").Value = 11217.88
$ws.Range("J132").Value = 7019.44
$ws.Range("K132").Value = 33653.64
$ws.Range("L132").Value = 21058.32
$ws.Range("M132").Value = -31123.64
$ws.Range("N132").Value = -26118.32
$ws.Range("H136").Value = 2266.392
$ws.Range("I136").Value = 1316.7317
$ws.Range("K136").Value = 3950.1951
$ws.Range("M136").Value = -1400.1951

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H80").Value = 42240.2
$ws.Range("J80").Value = 50000
$ws.Range("L80").Value = 50000
$ws.Range("N80").Value = -51996
$ws.Range("H83").Value = 42240.2
$ws.Range("J83").Value = 50000
$ws.Range("L83").Value = 150000
$ws.Range("N83").Value = -159984
$ws.Range("H107").Value = 816.3333
$ws.Range("J107").Value = 1150
$ws.Range("L107").Value = 3450
$ws.Range("N107").Value = -7290
$ws.Range("H109").Value = 27777
$ws.Range("J109").Value = 27777
$ws.Range("L109").Value = 27777
$ws.Range("N109").Value = -30551
$ws.Range("H110").Value = 39750
$ws.Range("J110").Value = 39750
$ws.Range("L110").Value = 39750
$ws.Range("N110").Value = -47930
$ws.Range("H111").Value = 29900
$ws.Range("J111").Value = 29900
$ws.Range("L111").Value = 29900
$ws.Range("N111").Value = -38080
$ws.Range("H114").Value = 39800
$ws.Range("J114").Value = 39800
$ws.Range("L114").Value = 39800
$ws.Range("N114").Value = -48478
$ws.Range("H116").Value = 35000
$ws.Range("J116").Value = 35000
$ws.Range("L116").Value = 35000
$ws.Range("N116").Value = -44178
$ws.Range("H117").Value = 0
$ws.Range("J117").Value = 0
$ws.Range("L117").Value = 0
$ws.Range("N117").ClearContents()
$ws.Range("H121").Value = 0
$ws.Range("J121").Value = 0
$ws.Range("L121").Value = 0
$ws.Range("N121").ClearContents()
$ws.Range("H132").Value = 6804357.5
$ws.Range("I132").Value = 761.4828
$ws.Range("J132").Value = 16669572
$ws.Range("K132").Value = 2284.4484
$ws.Range("L132").Value = 50008716
$ws.Range("M132").Value = 245.5515999999998
$ws.Range("N132").Value = -50013776
$ws.Range("H136").Value = 1745.4595
$ws.Range("I136").Value = 628.5102000000001
$ws.Range("J136").Value = 3934.68
$ws.Range("K136").Value = 1885.5306
$ws.Range("L136").Value = 11804.04
$ws.Range("M136").Value = 664.4694
$ws.Range("N136").Value = -16904.04
